$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2024-10-28 18:36:25"
$ws.Range("B6").Value = "ORM-0516041"
$ws.Range("C6").Value = "MCULE-2227031507"
$ws.Range("D6").Value = "C21H17CLN2O2"
$ws.Range("E6").Value = "Duplicate"
$ws.Range("F6").Value = "/home/robekott/ERAT/MCULE/mcule_test.sdf"

$ws.Range("A7").Value = "2024-10-28 18:36:27"
$ws.Range("B7").Value = "ORM-0516042"
$ws.Range("C7").Value = "MCULE-3988458386"
$ws.Range("D7").Value = "C16H14CLN5O"
$ws.Range("E7").Value = "Duplicate"
$ws.Range("F7").Value = "/home/robekott/ERAT/MCULE/mcule_test.sdf"
